$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 ("Save"), formatted like the other header cells (e.g. G1 "sum")
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New data cells for the "Save" column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
